$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '41.526.55'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +1.04%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.483.85'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '314.45'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.71%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '93.40'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.80%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.544'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -1.22%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.508'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +3.00%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '32.72'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -0.87%  '

$ws.Range("E11").Value = '  +1.63%  '

$ws.Range("E12").Value = '  +3.13%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '2.866.20'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +1.09%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '6.85'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -1.24%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '16.17'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +10.30%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '2.475.88'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.763'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -1.74%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '41.546.30'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '6.40'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +2.51%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '0.0₃0936'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +2.63%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '71.68'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +5.80%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '11.38'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +2.52%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '236.99'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +1.18%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.71'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("E27").Value = '  +5.23%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '2.20'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.22%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '9.68'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +1.59%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '35.97'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +1.37%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '157.73'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +4.97%  '

$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("E33").Value = '  +1.17%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.0753'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +2.64%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '17.45'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +4.11%  '

$ws.Range("E36").Value = '  -7.46%  '

$ws.Range("E37").Value = '  -0.87%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.105'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +3.04%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.84'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -0.69%  '

$ws.Range("E40").Value = '  +0.84%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '4.11'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +0.34%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '19.93'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -0.41%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '1.977.68'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +0.97%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.0284'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +1.40%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '2.95'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -1.39%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '9.08'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +7.20%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '2.723.17'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +1.18%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '98.05'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +2.27%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '68.04'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -2.75%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '72.45'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -1.60%  '
